# Update line-MW results for the 380 kV case (Case_2_189, pl_mw.xlsx, sheet "res_line")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per data row (row 2 = index 0 ... row 25 = index 23),
# only columns B, C, E, F, G, H, I, K, L, O change; others stay as-is.
$newValues = @{
    2 = @{ "B" = 0.4696084445338329; "C" = 0.1727505511875194; "E" = 0.113611068416791; "F" = 0.4443680307746121; "G" = 0.4738395392023698; "H" = 0.642644101790161; "I" = 0.6392682196701021; "K" = 0.3095076288651057; "L" = 0.1939107927522343; "O" = 2.19233016802594 }
    3 = @{ "B" = 0.4226914432631474; "C" = 0.1730467968707003; "E" = 0.1132844283242171; "F" = 0.3878228170618172; "G" = 0.4808712838512683; "H" = 0.6494570183898709; "I" = 0.6481425099264619; "K" = 0.2701769831902538; "L" = 0.1863756650207904; "O" = 2.221288052062476 }
    4 = @{ "B" = 0.3938636975041732; "C" = 0.1732658888122032; "E" = 0.1131468847916715; "F" = 0.3531389305169483; "G" = 0.4855348430984634; "H" = 0.6539160963363457; "I" = 0.6539444975516346; "K" = 0.2459323087983876; "L" = 0.1818436040563682; "O" = 2.24037239555706 }
    5 = @{ "B" = 0.3821119031680382; "C" = 0.1733645643709991; "E" = 0.1131067073818031; "F" = 0.3390132514313251; "G" = 0.4875222204166718; "H" = 0.6558026246420887; "I" = 0.6563976428044676; "K" = 0.2360291180359297; "L" = 0.1800206462387735; "O" = 2.248477253444761 }
    6 = @{ "B" = 0.3801602964326491; "C" = 0.1733815179314391; "E" = 0.113100995420865; "F" = 0.336668177824194; "G" = 0.4878574715926938; "H" = 0.6561200749448233; "I" = 0.6568103479660721; "K" = 0.2343833155256618; "L" = 0.1797193921832161; "O" = 2.249842856580493 }
    7 = @{ "B" = 0.393705224718417; "C" = 0.1732671814957527; "E" = 0.1131462786410324; "F" = 0.3529483938344953; "G" = 0.4855612936506901; "H" = 0.6539412575990475; "I" = 0.6539772220407407; "K" = 0.2457988442390047; "L" = 0.1818189221032185; "O" = 2.240480373240317 }
    8 = @{ "B" = 0.4534362488367663; "C" = 0.1728450008420737; "E" = 0.1134853800624533; "F" = 0.4248636149813478; "G" = 0.4761922194850783; "H" = 0.6449359502162082; "I" = 0.642254775425469; "K" = 0.2959666654931823; "L" = 0.1912931295607763; "O" = 2.202044131635205 }
    9 = @{ "B" = 0.5703717153257912; "C" = 0.1723105118881065; "E" = 0.1146493396939086; "F" = 0.5661985755041457; "G" = 0.4605684474345182; "H" = 0.629464133798507; "I" = 0.6220693341356132; "K" = 0.3935606243476286; "L" = 0.2106179166482605; "O" = 2.137021659034502 }
    10 = @{ "B" = 0.6561267087531917; "C" = 0.1720944533493451; "E" = 0.1158076001301076; "F" = 0.6702781546542269; "G" = 0.4507698776025109; "H" = 0.6194280066551343; "I" = 0.6089474072493193; "K" = 0.4647545103129005; "L" = 0.2252665297110497; "O" = 2.0955640023585 }
    11 = @{ "B" = 0.6950972629402656; "C" = 0.1720340518237933; "E" = 0.1164001207789731; "F" = 0.7176906081379002; "G" = 0.4466780286532313; "H" = 0.6151507795609703; "I" = 0.6033488374621072; "K" = 0.4970261652282488; "L" = 0.232027666177288; "O" = 2.078075828339934 }
    12 = @{ "B" = 0.709847864384642; "C" = 0.1720165880161204; "E" = 0.1166339055421837; "F" = 0.7356546913071611; "G" = 0.4451812030651467; "H" = 0.6135725216414158; "I" = 0.601282118639805; "K" = 0.5092294488666482; "L" = 0.2346018378991346; "O" = 2.071650785590407 }
    13 = @{ "B" = 0.7066713719980271; "C" = 0.1720201092108837; "E" = 0.1165831377173987; "F" = 0.7317853510981394; "G" = 0.445501227189304; "H" = 0.6139105851347324; "I" = 0.6017248506972823; "K" = 0.506602034888175; "L" = 0.2340468289138187; "O" = 2.073025753439651 }
    14 = @{ "B" = 0.6963109436856882; "C" = 0.1720325068844843; "E" = 0.1164191659521414; "F" = 0.7191683204515869; "G" = 0.4465538277781249; "H" = 0.6150201050002622; "I" = 0.6031777380725192; "K" = 0.4980304877467461; "L" = 0.2322391674904765; "O" = 2.077543279987651 }
    15 = @{ "B" = 0.6899639805436948; "C" = 0.1720408040751238; "E" = 0.1163199531932122; "F" = 0.7114413442032514; "G" = 0.4472054380287105; "H" = 0.6157051135654399; "I" = 0.6040746208402314; "K" = 0.4927778925304267; "L" = 0.2311337258453392; "O" = 2.080336101555147 }
    16 = @{ "B" = 0.653579006627524; "C" = 0.1720991597388348; "E" = 0.1157701957512636; "F" = 0.6671810134426437; "G" = 0.4510446543457363; "H" = 0.6197133326168327; "I" = 0.6093207488005614; "K" = 0.4626431046441439; "L" = 0.2248266224550264; "O" = 2.09673450185953 }
    17 = @{ "B" = 0.6312471093685303; "C" = 0.1721446344801265; "E" = 0.1154497261017369; "F" = 0.6400460337125793; "G" = 0.4534935920048397; "H" = 0.622246063447264; "I" = 0.6126340401910113; "K" = 0.4441264267430256; "L" = 0.2209822743866283; "O" = 2.107145697555737 }
    18 = @{ "B" = 0.6183987025596878; "C" = 0.1721743585268172; "E" = 0.1152715786915124; "F" = 0.6244449056556647; "G" = 0.4549365593706298; "H" = 0.6237299563308056; "I" = 0.6145746438616353; "K" = 0.4334653644099546; "L" = 0.2187802835227757; "O" = 2.113262987557079 }
    19 = @{ "B" = 0.6140478487541543; "C" = 0.1721850366035724; "E" = 0.1152123229600086; "F" = 0.619163680173358; "G" = 0.4554310290386496; "H" = 0.6242370375891113; "I" = 0.615237689495391; "K" = 0.4298538879450007; "L" = 0.2180363073074716; "O" = 2.115356354935116 }
    20 = @{ "B" = 0.6336247669457578; "C" = 0.1721394246053407; "E" = 0.115483201407045; "F" = 0.642933953830422; "G" = 0.4532293370794136; "H" = 0.621973641868621; "I" = 0.6122777238141222; "K" = 0.4460986779105838; "L" = 0.2213905630208046; "O" = 2.10602405093104 }
    21 = @{ "B" = 0.6993542422610801; "C" = 0.1720287188980834; "E" = 0.1164670732590807; "F" = 0.7228739723491628; "G" = 0.4462432232579232; "H" = 0.614693087708801; "I" = 0.6027495419544735; "K" = 0.5005486344408325; "L" = 0.232769745924486; "O" = 2.076211015477469 }
    22 = @{ "B" = 0.7422727410747711; "C" = 0.171987875953441; "E" = 0.1171649270205286; "F" = 0.7751780083420101; "G" = 0.4419844340964403; "H" = 0.6101763332122161; "I" = 0.5968332365663969; "K" = 0.5360336778390717; "L" = 0.2402875188850686; "O" = 2.057876859764264 }
    23 = @{ "B" = 0.7193702674868518; "C" = 0.1720068042127281; "E" = 0.1167874605303965; "F" = 0.7472568307830727; "G" = 0.4442293028961046; "H" = 0.6125649172593839; "I" = 0.5999624192551938; "K" = 0.5171041558074876; "L" = 0.2362677905081938; "O" = 2.067556823005006 }
    24 = @{ "B" = 0.6325498574775565; "C" = 0.1721417688347415; "E" = 0.1154680482349804; "F" = 0.6416283278902171; "G" = 0.4533486975910108; "H" = 0.6220967170923259; "I" = 0.6124387030982934; "K" = 0.4452070716206151; "L" = 0.2212059502524255; "O" = 2.106530736979352 }
    25 = @{ "B" = 0.5387625360448851; "C" = 0.1724239267075944; "E" = 0.1142811506957671; "F" = 0.5279251897347166; "G" = 0.464500342150636; "H" = 0.6334157440966806; "I" = 0.6272300319281197; "K" = 0.3672460074113815; "L" = 0.2053106109695051; "O" = 2.153503392370155 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range("$colLetter$rowNum").Value = [double]$rowData[$colLetter]
    }
}

